$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Status and Date values ---
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Cells.Item(6, 2).Value = "draft"
$ws1.Cells.Item(8, 2).Value = "2023-08-01T16:12:28+00:00"

# --- Sheet "Include from duo.owl": rename + restructure to hold Ferlab.bio codes ---
$ws2 = $wb.Worksheets.Item("Include from duo.owl")
$ws2.Name = "Include from Ferlab.bio CodeS"

# Remove the four DUO concept detail rows (old rows 3-6), shifting following rows up
$ws2.Range("A3:B6").EntireRow.Delete()

# Row 1 header: "Concept" -> "Codes"; drop the "Description" header in B1
$ws2.Cells.Item(1, 1).Value = "Codes"
$ws2.Cells.Item(1, 2).Clear()

# Row 2: "DUO:0000004" -> "All codes"; drop the "No restriction" value in B2
$ws2.Cells.Item(2, 1).Value = "All codes"
$ws2.Cells.Item(2, 2).Clear()
